$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates -------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 2.0.0-ballot -> 2.0.0
$wsMeta.Cells.Item(3, 2).Value = "2.0.0"

# Date: refreshed publication date/time
$wsMeta.Cells.Item(8, 2).Value = "2025-10-20T13:10:23+00:00"

# --- Elements sheet updates --------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

# Remove the now-obsolete child element rows (Extension.value[x].id/.extension/
# .coding/.text) that described the previous CodeableConcept sub-structure.
$wsElem.Rows("7:10").Delete()

# Row 6 (Extension.value[x]) gains a binding strength + value set, matching the
# newly introduced required binding to the "motif non realisation" value set.
$wsElem.Range("X6").Value = "required"
$wsElem.Range("Y6").Value = ""
$wsElem.Range("Z6").Value = "https://smt.esante.gouv.fr/fhir/ValueSet/jdv-motif-non-realisation-evenement-cisis"
